$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERT")
Write-Host ($ws.Columns.Item(6) | Get-Member | Out-String)
